$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A25").Value = 2016
$ws.Range("B25").Value = "Previsión población 2016"
$ws.Range("B25").Font.Name = "Arial"
$ws.Range("B25").Font.Size = 10
$ws.Range("B25").Font.Bold = $false
$ws.Range("C25").Formula = "=B21*C23+B21"
$ws.Range("C25").NumberFormat = "#,##0"
$ws.Range("C25").Select()
